$d = $word.ActiveDocument

# 1. Remove the existing "_GoBack" bookmark from the first paragraph
#    (it currently sits between the two <w:tab/> runs there).
$goBack = $d.Bookmarks("_GoBack")
$goBack.Delete()

# 2. Append a new paragraph at the end of the body (after the "Assigned" /
#    tab-stops paragraph, before the sectPr) that demonstrates negative /
#    positive tab-stop positions together with a hanging indent, and that
#    re-creates the "_GoBack" bookmark in its new location.
$lastPara = $d.Paragraphs.Last
$endOfBody = $d.Range($lastPara.Range.End, $lastPara.Range.End)

$newParagraphXml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
  <w:pPr>
    <w:tabs>
      <w:tab w:val="left" w:pos="-720"/>
      <w:tab w:val="left" w:pos="720"/>
    </w:tabs>
    <w:ind w:hanging="1080"/>
  </w:pPr>
  <w:r>
    <w:tab/>
  </w:r>
  <w:r>
    <w:t>-0.5</w:t>
  </w:r>
  <w:r>
    <w:tab/>
  </w:r>
  <w:r>
    <w:tab/>
  </w:r>
  <w:bookmarkStart w:id="0" w:name="_GoBack"/>
  <w:bookmarkEnd w:id="0"/>
  <w:r>
    <w:t>0.5</w:t>
  </w:r>
</w:p>
'@

[void]$endOfBody.InsertXML($newParagraphXml)
